# Updates cryptos price/volume data per upstream source refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.349.01'
$ws.Range("E2").Value = '  +1.84%  '
$ws.Range("D3").Value = '1.865.60'
$ws.Range("E3").Value = '  +1.44%  '
$ws.Range("E4").Value = '  +1.58%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.21'
$ws.Range("E5").Value = '  +1.68%  '
$ws.Range("E6").Value = '  +1.40%  '
$ws.Range("E7").Value = '  +1.95%  '
$ws.Range("E8").Value = '  +2.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07407'
$ws.Range("E9").Value = '  +3.69%  '
$ws.Range("E10").Value = '  +2.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.79'
$ws.Range("E11").Value = '  +6.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07890'
$ws.Range("E12").Value = '  +3.90%  '
$ws.Range("D13").Value = '1.884.19'
$ws.Range("E13").Value = '  +1.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.440'
$ws.Range("E14").Value = '  +3.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.559'
$ws.Range("E15").Value = '  +2.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.32'
$ws.Range("E16").Value = '  +3.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.022'
$ws.Range("E17").Value = '  +1.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008799'
$ws.Range("E18").Value = '  +2.09%  '
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.82'
$ws.Range("E20").Value = '  +2.41%  '
$ws.Range("D21").Value = '27.399.19'
$ws.Range("E21").Value = '  +1.87%  '
$ws.Range("E22").Value = '  +2.65%  '
$ws.Range("E23").Value = '  +1.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.966'
$ws.Range("E24").Value = '  +2.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.29'
$ws.Range("E25").Value = '  +1.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.55'
$ws.Range("E26").Value = '  +2.19%  '
$ws.Range("E27").Value = '  +0.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '116.12'
$ws.Range("E28").Value = '  +1.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.999'
$ws.Range("E29").Value = '  +3.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08933'
$ws.Range("E30").Value = '  +1.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.354'
$ws.Range("E31").Value = '  +3.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.197'
$ws.Range("E32").Value = '  +2.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.569'
$ws.Range("E33").Value = '  +2.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7447'
$ws.Range("E34").Value = '  +0.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.693'
$ws.Range("E35").Value = '  -1.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02053'
$ws.Range("E36").Value = '  +5.78%  '
$ws.Range("E37").Value = '  +3.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05303'
$ws.Range("E38").Value = '  +1.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5389'
$ws.Range("E39").Value = '  +4.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.120'
$ws.Range("E40").Value = '  +2.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1528'
$ws.Range("E41").Value = '  +1.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.405'
$ws.Range("E42").Value = '  +3.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.71'
$ws.Range("E43").Value = '  +2.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4839'
$ws.Range("E44").Value = '  +3.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.021'
$ws.Range("E45").Value = '  +1.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.674'
$ws.Range("E46").Value = '  +5.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.43'
$ws.Range("E47").Value = '  +1.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '66.91'
$ws.Range("E48").Value = '  +2.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06099'
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9017'
$ws.Range("E50").Value = '  +2.13%  '
$ws.Range("E51").Value = '  +1.93%  '
